$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

# Copy AC27's existing format down to AC28:AC62, then set all the values
$src = $ws.Range("AC27")
$dst = $ws.Range("AC28:AC62")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("AC27:AC62").Value = "UIResources/BuilderSprite"
